# Refresh Market Board snapshot columns (H:N) across the per-job leve-profit
# sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW) with the latest scheduled-runner pull.
# Columns: H=currentAveragePrice, I/J=NQ/HQ average price, K/L=Leve price NQ/HQ,
# M/N=Leve profit NQ/HQ. A cleared cell means that price/profit no longer applies.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 1915
$ws.Range("I58").Value = 80
$ws.Range("J58").Value = 3750
$ws.Range("K58").Value = 240
$ws.Range("L58").Value = 11250
$ws.Range("M58").Value = -90
$ws.Range("N58").Value = -11550
# Row 95
$ws.Range("H95").Value = 34000
$ws.Range("J95").Value = 34000
$ws.Range("L95").Value = 34000
$ws.Range("N95").Value = -39492
# Row 98
$ws.Range("H98").Value = 2288
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
# Row 122
$ws.Range("H122").Value = 2288
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
# Row 127
$ws.Range("H127").Value = 929.5
$ws.Range("I127").Value = 929.5
$ws.Range("K127").Value = 2788.5
$ws.Range("M127").Value = 2171.5
$ws = $wb.Worksheets.Item("ARM")
# Row 33
$ws.Range("H33").Value = 20000
$ws.Range("I33").Value = 20000
$ws.Range("K33").Value = 20000
$ws.Range("M33").Value = -19671
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
# Row 61
$ws.Range("H61").Value = 1499.5
$ws.Range("I61").Value = 999
$ws.Range("K61").Value = 999
$ws.Range("M61").Value = -787
# Row 97
$ws.Range("H97").Value = 41669410
$ws.Range("I97").Value = 47621964
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 47621964
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -47621468
$ws.Range("N97").Value = -2492
# Row 122
$ws.Range("H122").Value = 7738.875
$ws.Range("I122").Value = 6985.1665
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 20955.4995
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -18505.4995
$ws.Range("N122").Value = -34900
# Row 132
$ws.Range("H132").Value = 9000
$ws.Range("I132").Value = 9000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 27000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -24470
$ws.Range("N132").ClearContents()
# Row 136
$ws.Range("H136").Value = 1499.5
$ws.Range("I136").Value = 999
$ws.Range("K136").Value = 2997
$ws.Range("M136").Value = -447
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1574.7778
$ws.Range("I134").Value = 1310.4286
$ws.Range("K134").Value = 3931.2858
$ws.Range("M134").Value = -1396.2858
$ws = $wb.Worksheets.Item("CRP")
# Row 64
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
# Row 67
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 141
$ws.Range("H141").Value = 705554.6
$ws.Range("J141").Value = 705554.6
$ws.Range("L141").Value = 705554.6
$ws.Range("N141").Value = -715914.6
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 459.63635
$ws.Range("I2").Value = 595.25
$ws.Range("J2").Value = 98
$ws.Range("K2").Value = 3571.5
$ws.Range("L2").Value = 588
$ws.Range("M2").Value = -3458.5
$ws.Range("N2").Value = -814
# Row 22
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3338
# Row 27
$ws.Range("H27").Value = 1000
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3204
# Row 97
$ws.Range("H97").Value = 1241.5
$ws.Range("I97").Value = 999.6667
$ws.Range("K97").Value = 2999.0001
$ws.Range("M97").Value = -2503.0001
# Row 98
$ws.Range("H98").Value = 3313
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 3313
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 9939
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -12935
# Row 122
$ws.Range("H122").Value = 606.25
$ws.Range("I122").Value = 606.25
$ws.Range("K122").Value = 5456.25
$ws.Range("M122").Value = -3006.25
# Row 131
$ws.Range("H131").Value = 1162
$ws.Range("I131").Value = 855.6667
$ws.Range("K131").Value = 2567.0001
$ws.Range("M131").Value = 2472.9999
$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 7703.2
$ws.Range("I22").Value = 1505.3334
$ws.Range("J22").Value = 17000
$ws.Range("K22").Value = 1505.3334
$ws.Range("L22").Value = 17000
$ws.Range("M22").Value = -976.3334
$ws.Range("N22").Value = -18058
# Row 39
$ws.Range("H39").Value = 43000
$ws.Range("J39").Value = 43000
$ws.Range("L39").Value = 43000
$ws.Range("N39").Value = -44064
# Row 43
$ws.Range("H43").Value = 16602.428
$ws.Range("I43").Value = 2804.25
$ws.Range("J43").Value = 35000
$ws.Range("K43").Value = 2804.25
$ws.Range("L43").Value = 35000
$ws.Range("M43").Value = -2653.25
$ws.Range("N43").Value = -35302
# Row 46
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2844
# Row 128
$ws.Range("H128").Value = 89999
$ws.Range("J128").Value = 89999
$ws.Range("L128").Value = 89999
$ws.Range("N128").Value = -99959
# Row 132
$ws.Range("H132").Value = 2782.6
$ws.Range("I132").Value = 2782.6
$ws.Range("K132").Value = 8347.799999999999
$ws.Range("M132").Value = -5817.799999999999
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 461.7097
$ws.Range("I22").Value = 227.66667
$ws.Range("K22").Value = 227.66667
$ws.Range("M22").Value = 67.33332999999999
# Row 27
$ws.Range("H27").Value = 461.7097
$ws.Range("I27").Value = 227.66667
$ws.Range("K27").Value = 227.66667
$ws.Range("M27").Value = -120.66667
# Row 38
$ws.Range("H38").Value = 61500
$ws.Range("I38").Value = 90000
$ws.Range("J38").Value = 33000
$ws.Range("K38").Value = 90000
$ws.Range("L38").Value = 33000
$ws.Range("M38").Value = -89590
$ws.Range("N38").Value = -33820
# Row 39
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 5000
$ws.Range("K39").Value = 5000
$ws.Range("M39").Value = -4540
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 47
$ws.Range("H47").Value = 23500
$ws.Range("J47").Value = 37000
$ws.Range("L47").Value = 37000
$ws.Range("N47").Value = -37980
# Row 52
$ws.Range("H52").Value = 23500
$ws.Range("J52").Value = 37000
$ws.Range("L52").Value = 37000
$ws.Range("N52").Value = -37466
